$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Task Distribution"
$ws.Range("C4").Value = "9/13/2024 9:30pm - 10:00pm"
$ws.Range("E4").Value = "Sayantika, Manisha"

$ws.Range("E5:F5").Select()
